# Update the "previsao_retorno" workbook for BIBI/PF with refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: nome_cliente -> nome (column K)
$ws.Range("K1").Value = "nome"

# Refresh "situacao" (column J) text for clients whose inactivity counter
# ticked forward since the previous export.
$ws.Range("J72").Value = "INATIVO - 31.9 meses sem comprar"
$ws.Range("J73").Value = "INATIVO - 6.2 meses sem comprar"
$ws.Range("J80").Value = "INATIVO - 24.7 meses sem comprar"
$ws.Range("J85").Value = "INATIVO - 3.4 meses sem comprar"
$ws.Range("J89").Value = "INATIVO - 31.9 meses sem comprar"
$ws.Range("J90").Value = "INATIVO - 12.2 meses sem comprar"
$ws.Range("J92").Value = "INATIVO - 14.9 meses sem comprar"
$ws.Range("J96").Value = "INATIVO - 0.7 meses sem comprar"

# Row 110 (BEMOL S/A) - updated purchase history count and refreshed
# ultima_compra / proxima_compra timestamps.
$ws.Range("E110").Value = 14254
$ws.Range("H110").Value = 45789.75063657408
$ws.Range("I110").Value = 45790.75063657408
